$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4739503264427185
$ws.Range("B1").Value = 0.4614225625991821
$ws.Range("C1").Value = 0.4774905443191528
$ws.Range("D1").Value = 0.637269139289856
$ws.Range("E1").Value = 0.7532960772514343
